$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'255.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.02%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.91%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.233"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.47%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05863"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.08%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.715"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.92%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8688"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.30%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.76%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.45%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.85%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03182"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.23%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09235"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.46%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001552"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.95%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006051"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-94.06%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005808"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.36%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.506"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.23%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.229"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.34%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.222"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.17%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3175"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.79%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03458"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'3.55%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.25%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.526"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.81%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04163"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.09%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1375"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.49%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004796"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'15.33%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.09%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'1.19%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03815"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.84%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.005666"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.07%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002355"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'12.13%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01045"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.79%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.25%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'11.24%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-13.62%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
